# Applies the tile/map key updates described in the commit message:
# "changed a bunch of the tiles and added a bunch. Also modified the maps
#  to reflect the changes. Updated the wbs a little, and added the key to
#  what each character means in the maps"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New "key" values filled into the previously-empty Map/Key column (G)
# for several tile rows - these introduce the new shared strings:
# "done", "done-late", "late 1month", "needs pathfinding"
$ws.Range("G28").Value = "done"
$ws.Range("G29").Value = "done"
$ws.Range("G30").Value = "done"
$ws.Range("G31").Value = "done"

$ws.Range("G34").Value = "done"
$ws.Range("G35").Value = "done"

$ws.Range("G38").Value = "done-late"
$ws.Range("G39").Value = "late 1month"

$ws.Range("G41").Value = "done"
$ws.Range("G42").Value = "done"

$ws.Range("G44").Value = "done"
$ws.Range("G45").Value = "needs pathfinding"
$ws.Range("G46").Value = "done"

$ws.Range("G49").Value = "done"
$ws.Range("G50").Value = "done"
$ws.Range("G51").Value = "In Progress"

$ws.Range("G54").Value = "done"
$ws.Range("G55").Value = "done"

# Row 56 task (Trade Window) moved from "Not Started"/"In Progress" to "Done"
$ws.Range("B56").Value = "Done"
$ws.Range("G56").Value = "Done"

# Update the selected cell on Sheet1 to reflect where the author last worked
$ws.Range("B51").Select()
